# "9th Stab - Cosmetic Changes"
# Two new rating-date columns (Jun_17, Jun_15) are inserted right after the
# firm-name column (B), pushing the existing Jun_13 / Jun_10 columns from
# B/C to D/E. The new columns are filled with the same "UN" (unrated)
# placeholder used throughout the sheet, except for Zacks Investment
# Research (row 5) which records a real upgrade on 6/15/2018, highlighted
# with a fill color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing date columns (old B:C -> new D:E)
$ws.Columns("B:C").Insert()

# New date headers
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the two new columns with the "UN" (unrated) placeholder for every
# analyst row
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Zacks Investment Research (row 5) actually moved from Sell to Hold on
# 6/15/2018 - record it in the new Jun_15 column with a highlight fill
$ws.Range("C5").Value = "6/15/2018,Upgrades,Sell -> Hold,"
$ws.Range("C5").Interior.ColorIndex = 42

# Match the original column widths (8 characters) for the date columns
$ws.Columns("C").ColumnWidth = 7.1666666666666666
$ws.Columns("D").ColumnWidth = 7.1666666666666666
$ws.Columns("E").ColumnWidth = 7.1666666666666666
